$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Type")
$v = $ws.Range("C11").Value
Write-Host "Value: $v"
Write-Host "Type: $($v.GetType())"
